$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.281.52"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.273.92"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.92"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.78"
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.87"
$ws.Range("E10").Value = "  +10.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.68"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.599.47"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.46"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.270.49"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.170.63"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.89"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.93"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.59"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.80"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.50"
$ws.Range("E28").Value = "  +7.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.34"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.34"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.106"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.84"
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("E42").Value = "  +14.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.989.26"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.01"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.34"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.16"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.75"
$ws.Range("E51").Value = "  -0.01%  "
